$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Complete row 196 (previously only had A196 and G196 populated)
$ws.Range("B196").Value = 125000
$ws.Range("C196").Value = 1.5
$ws.Range("D196").Value = 1.5
$ws.Range("E196").Value = 1.5
$ws.Range("F196").Value = 6
$ws.Range("G196").Value = 1.5

# New row 197
$ws.Range("A197").Value = "13-10-2021"
$ws.Range("B197").Value = 365000
$ws.Range("C197").Value = 1.5
$ws.Range("D197").Value = 1.5
$ws.Range("E197").Value = 1.5
$ws.Range("F197").Value = 7
$ws.Range("G197").Value = 1.5

# New row 198
$ws.Range("A198").Value = "14-10-2021"
$ws.Range("B198").Value = 97000
$ws.Range("C198").Value = 2.75
$ws.Range("D198").Value = 2.75
$ws.Range("E198").Value = 2.75
$ws.Range("F198").Value = 5
$ws.Range("G198").Value = 2.75

# New row 199
$ws.Range("A199").Value = "15-10-2021"
$ws.Range("B199").Value = 97000
$ws.Range("C199").Value = 2.75
$ws.Range("D199").Value = 2.75
$ws.Range("E199").Value = 2.75
$ws.Range("F199").Value = 5
$ws.Range("G199").Value = 2.75

# New row 200 (only A and G populated)
$ws.Range("A200").Value = "18-10-2021"
$ws.Range("G200").Value = 2.75
